# Weekly fruit/vegetable price update: insert a new observation row at row 4
# (pushing the existing rows 4-75 down to 5-76) and populate it with the
# latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 4; this shifts existing rows 4..75 down
# to 5..76 and extends the used range accordingly.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with this week's record.
$ws.Cells.Item(4, 1).Value = 6
$ws.Cells.Item(4, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(4, 3).Value = "Metropolitana"
$ws.Cells.Item(4, 4).Value = 45083
$ws.Cells.Item(4, 5).Value = 13
$ws.Cells.Item(4, 6).Value = 100112035
$ws.Cells.Item(4, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 18000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 18800
$ws.Cells.Item(4, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(4, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(4, 16).Value = 1253
$ws.Cells.Item(4, 17).Value = 15
$ws.Cells.Item(4, 18).Value = "Hortaliza"
